# ---------------------------------------------------------------------------
# kep.xlsx edit: add a "variables" lookup sheet (code/description/unit) after
# "month", clear a stray placeholder cell on "month", and populate the new
# sheet's shared-string-backed text columns in the same column-major order the
# original authoring tool used (index column, then code, then description, then
# unit) so the shared-string table grows identically.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Clear the old placeholder value in month!AV3 (was referencing the
#        leftover "###" shared string; the column has no data for this row).
$monthSheet = $wb.Worksheets.Item("month")
$monthSheet.Cells.Item(3, 48).Value = $null

# --- 2. Add the new "variables" worksheet right after "month".
$ws = $wb.Worksheets.Add($null, $monthSheet)
$ws.Name = "variables"

# --- 3. Column data (values) for the lookup table, rows 2..50.
$colA = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48)
$colB = @('CONSTR_bln_rub_fix','CONSTR_rog','CONSTR_yoy','CORP_DEBT_bln_rub','CORP_DEBT_overdue','CORP_DEBT_rog','CPI_ALCOHOL_rog','CPI_FOOD_rog','CPI_NONFOOD_rog','CPI_SERVICES_rog','CPI_rog','IND_PROD_rog','IND_PROD_yoy','IND_PROD_ytd','I_bln_rub','I_rog','I_yoy','PROD_AUTO_BUS_units','PROD_AUTO_PSGR_th','PROD_AUTO_TRUCKS_th','PROD_BYCYCLES_th','PROD_E_TWh','PROD_RAILWAY_CARGO_WAGONS_units','PROD_RAILWAY_PSGR_WAGONS_units','RETAIL_SALES_bln_rub','RETAIL_SALES_rog','RETAIL_SALES_yoy','RUR_EUR_eop','RUR_USD_eop','SOC_EMPLOYED_mln','SOC_EMPLOYED_yoy','SOC_PENSION_rub','SOC_UNEMPLOYED_bln','SOC_UNEMPLOYMENT_percent','SOC_WAGE_rog','SOC_WAGE_rub','SOC_WAGE_yoy','TRANS_COM_bln_t_km','TRANS_COM_rog','TRANS_COM_yoy','TRANS_RAILLOAD_mln_t','TRANS_RAILLOAD_rog','TRANS_RAILLOAD_yoy','TRANS_bln_t_km','TRANS_rog','TRANS_yoy','USLUGI_bln_rub','USLUGI_rog','USLUGI_yoy')
$colC = @('Объем работ по виду деятельности "Строительство"','Объем работ по виду деятельности "Строительство"','Объем работ по виду деятельности "Строительство"','Кредиторская задолженность','Кредиторская задолженность','Кредиторская задолженность','алкогольные напитки','продукты питания','<...>','<...>','Индекс потребительских цен','Индекс промышленного производства','Индекс промышленного производства','Индекс промышленного производства','Инвестиции в основной капитал','Инвестиции в основной капитал','Инвестиции в основной капитал','Автобусы, штук','Автомобили легковые, тыс.штук','Грузовые автомобили, тыс.штук','Велосипеды (без детских), тыс.штук','Электроэнергия, млрд. кВт·ч','Вагоны грузовые магистральные, штук','Вагоны пассажирские магистральные, штук','Оборот розничной торговли','Оборот розничной торговли','Оборот розничной торговли','Официальный курс евро по отношению к рублю','Официальный курс доллара США','Численность занятого в экономике населения','Численность занятого в экономике населения','Средний размер назначенных пенсий','Общая численность безработных','Уровень безработицы','Среднемесячная номинальная начисленная заработная плата одного работника','Среднемесячная номинальная начисленная заработная плата одного работника','Среднемесячная номинальная начисленная заработная плата одного работника','Коммерческий грузооборот транспорта','Коммерческий грузооборот транспорта','Коммерческий грузооборот транспорта','Freight loading on railway transport','Freight loading on railway transport','Freight loading on railway transport','Грузооборот транспорта, включая коммерческий и некоммерческий грузооборот','Грузооборот транспорта, включая коммерческий и некоммерческий грузооборот','Грузооборот транспорта, включая коммерческий и некоммерческий грузооборот','Объем платных услуг населению','Объем платных услуг населению','Объем платных услуг населению')
$colD = @('млрд. руб. (в фикс. ценах)','в % к предыдущему периоду','в % к аналог. периоду предыдущего года','млрд. руб.','<...>','в % к предыдущему периоду','в % к предыдущему периоду','в % к предыдущему периоду','в % к предыдущему периоду','в % к предыдущему периоду','в % к предыдущему периоду','в % к предыдущему периоду','в % к аналог. периоду предыдущего года','<...>','млрд. руб.','в % к предыдущему периоду','в % к аналог. периоду предыдущего года','штук','тыс.','тыс.','тыс.','млрд. кВт·ч','штук','штук','млрд. руб.','в % к предыдущему периоду','в % к аналог. периоду предыдущего года','на конец периода','на конец периода','млн. человек','в % к аналог. периоду предыдущего года','рублей','млрд.','%','в % к предыдущему периоду','рублей','в % к аналог. периоду предыдущего года','млрд. т-км','в % к предыдущему периоду','в % к аналог. периоду предыдущего года','млн. т','в % к предыдущему периоду','в % к аналог. периоду предыдущего года','млрд. т-км','в % к предыдущему периоду','в % к аналог. периоду предыдущего года','млрд. руб.','в % к предыдущему периоду','в % к аналог. периоду предыдущего года')
$header = @('Код','Описание','Ед.изм.')

# --- 4. Header row (B1:D1).
$ws.Cells.Item(1, 2).Value = $header[0]
$ws.Cells.Item(1, 3).Value = $header[1]
$ws.Cells.Item(1, 4).Value = $header[2]

# --- 5. Write column-major: all of A, then all of B, then all of C, then all
#        of D -- this is the order new shared strings get interned in, and it
#        has to match so the shared string table indices line up.
for ($i = 0; $i -lt $colA.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $colC.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt $colD.Count; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}

# --- 6. Formatting: header row + index column are bold, bordered, centered
#        (mirrors the look used for header rows on the other sheets).
$lastRow = $colA.Count + 1
$hdr = $ws.Range("B1:D1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

$idxCol = $ws.Range("A2:A" + $lastRow)
$idxCol.Font.Bold = $true
$idxCol.Borders.LineStyle = 1
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160

